$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8497.25
$ws.Range("I32").Value = 14999.5
$ws.Range("J32").Value = 1995
$ws.Range("K32").Value = 14999.5
$ws.Range("L32").Value = 1995
$ws.Range("M32").Value = -14673.5
$ws.Range("N32").Value = -2647
$ws.Range("H41").Value = 526.1177
$ws.Range("J41").Value = 455.7143
$ws.Range("L41").Value = 455.7143
$ws.Range("N41").Value = -1335.7143
$ws.Range("H62").Value = 6526.923
$ws.Range("I62").Value = 6735
$ws.Range("K62").Value = 6735
$ws.Range("M62").Value = -6111
$ws.Range("H65").Value = 6526.923
$ws.Range("I65").Value = 6735
$ws.Range("K65").Value = 33675
$ws.Range("M65").Value = -30555
$ws.Range("H107").Value = 1538.3572
$ws.Range("I107").Value = 1044.4166
$ws.Range("J107").Value = 4502
$ws.Range("K107").Value = 1044.4166
$ws.Range("L107").Value = 4502
$ws.Range("M107").Value = 875.5834
$ws.Range("N107").Value = -8342
$ws.Range("H129").Value = 4852
$ws.Range("I129").Value = 675.1667
$ws.Range("K129").Value = 2025.5001
$ws.Range("M129").Value = 2974.4999
$ws.Range("H133").Value = 51156
$ws.Range("J133").Value = 51156
$ws.Range("L133").Value = 51156
$ws.Range("N133").Value = -61276
$ws.Range("H138").Value = 2474.0205
$ws.Range("J138").Value = 2587.3845
$ws.Range("L138").Value = 7762.1535
$ws.Range("N138").Value = -18042.1535

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1950
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3340
$ws.Range("H32").Value = 4780.8857
$ws.Range("I32").Value = 2897.5386
$ws.Range("K32").Value = 2897.5386
$ws.Range("M32").Value = -2610.5386
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H98").Value = 24950
$ws.Range("J98").Value = 24950
$ws.Range("L98").Value = 24950
$ws.Range("N98").Value = -30940
$ws.Range("H122").Value = 1528.4
$ws.Range("I122").Value = 1528.4
$ws.Range("K122").Value = 4585.200000000001
$ws.Range("M122").Value = -2135.200000000001
$ws.Range("M51").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26700.65
$ws.Range("I82").Value = 12242.25
$ws.Range("J82").Value = 30315.25
$ws.Range("K82").Value = 12242.25
$ws.Range("L82").Value = 30315.25
$ws.Range("M82").Value = -11859.25
$ws.Range("N82").Value = -31081.25
$ws.Range("H85").Value = 26700.65
$ws.Range("I85").Value = 12242.25
$ws.Range("J85").Value = 30315.25
$ws.Range("K85").Value = 12242.25
$ws.Range("L85").Value = 30315.25
$ws.Range("M85").Value = -10916.25
$ws.Range("N85").Value = -32967.25
$ws.Range("H94").Value = 944.65955
$ws.Range("I94").Value = 944.65955
$ws.Range("K94").Value = 944.65955
$ws.Range("M94").Value = -493.65955
$ws.Range("H107").Value = 2583.3333
$ws.Range("I107").Value = 2461.3333
$ws.Range("K107").Value = 2461.3333
$ws.Range("M107").Value = -541.3332999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2980.8696
$ws.Range("I31").Value = 1390.2
$ws.Range("J31").Value = 3174.8538
$ws.Range("K31").Value = 1390.2
$ws.Range("L31").Value = 3174.8538
$ws.Range("M31").Value = -1095.2
$ws.Range("N31").Value = -3764.8538
$ws.Range("H34").Value = 2980.8696
$ws.Range("I34").Value = 1390.2
$ws.Range("J34").Value = 3174.8538
$ws.Range("K34").Value = 1390.2
$ws.Range("L34").Value = 3174.8538
$ws.Range("M34").Value = -1188.2
$ws.Range("N34").Value = -3578.8538
$ws.Range("H58").Value = 4456.6772
$ws.Range("I58").Value = 1505
$ws.Range("K58").Value = 1505
$ws.Range("M58").Value = -1302
$ws.Range("H86").Value = 3154.6667
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H89").Value = 3154.6667
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("H94").Value = 1365.7
$ws.Range("I94").Value = 1521.4
$ws.Range("J94").Value = 1313.8
$ws.Range("K94").Value = 1521.4
$ws.Range("L94").Value = 1313.8
$ws.Range("M94").Value = -1070.4
$ws.Range("N94").Value = -2215.8
$ws.Range("H99").Value = 12090.059
$ws.Range("I99").Value = 11013.091
$ws.Range("K99").Value = 11013.091
$ws.Range("M99").Value = -9515.091
$ws.Range("H107").Value = 5013
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5013
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5013
$ws.Range("N107").Value = -8853
$ws.Range("I122").Value = 1748.5
$ws.Range("J122").Value = 1397
$ws.Range("K122").Value = 5245.5
$ws.Range("L122").Value = 4191
$ws.Range("M122").Value = -2795.5
$ws.Range("N122").Value = -9091
$ws.Range("H126").Value = 12090.059
$ws.Range("I126").Value = 11013.091
$ws.Range("K126").Value = 33039.273
$ws.Range("M126").Value = -30569.273
$ws.Range("H136").Value = 4456.6772
$ws.Range("I136").Value = 1505
$ws.Range("K136").Value = 4515
$ws.Range("M136").Value = -1965
$ws.Range("M107").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2823.5
$ws.Range("I8").Value = 2823.5
$ws.Range("K8").Value = 8470.5
$ws.Range("M8").Value = -8331.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1074666.4
$ws.Range("J11").Value = 1371428.6
$ws.Range("L11").Value = 1371428.6
$ws.Range("N11").Value = -1371706.6
$ws.Range("H12").Value = 17499
$ws.Range("I12").Value = 14998
$ws.Range("J12").Value = 20000
$ws.Range("K12").Value = 14998
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = -14858
$ws.Range("N12").Value = -20280
$ws.Range("H80").Value = 5781.0713
$ws.Range("I80").Value = 2327.9167
$ws.Range("J80").Value = 26500
$ws.Range("K80").Value = 2327.9167
$ws.Range("L80").Value = 26500
$ws.Range("M80").Value = -1329.9167
$ws.Range("N80").Value = -28496
$ws.Range("H83").Value = 5781.0713
$ws.Range("I83").Value = 2327.9167
$ws.Range("J83").Value = 26500
$ws.Range("K83").Value = 11639.5835
$ws.Range("L83").Value = 132500
$ws.Range("M83").Value = -6647.583500000001
$ws.Range("N83").Value = -142484
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("H102").Value = 3010.875
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("H126").Value = 3986.75
$ws.Range("I126").Value = 3489
$ws.Range("K126").Value = 10467
$ws.Range("M126").Value = -7997
$ws.Range("H139").Value = 73858.7
$ws.Range("J139").Value = 73858.7
$ws.Range("L139").Value = 73858.7
$ws.Range("N139").Value = -84138.7
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1212.5641
$ws.Range("J46").Value = 1366.8
$ws.Range("L46").Value = 1366.8
$ws.Range("N46").Value = -1742.8
$ws.Range("H61").Value = 2119.6667
$ws.Range("I61").Value = 2286.5557
$ws.Range("K61").Value = 2286.5557
$ws.Range("M61").Value = -2084.5557
$ws.Range("H68").Value = 2732.5518
$ws.Range("I68").Value = 2119.3044
$ws.Range("K68").Value = 2119.3044
$ws.Range("M68").Value = -1370.3044
$ws.Range("H71").Value = 2732.5518
$ws.Range("I71").Value = 2119.3044
$ws.Range("K71").Value = 10596.522
$ws.Range("M71").Value = -6852.522000000001
$ws.Range("H113").Value = 2119.6667
$ws.Range("I113").Value = 2286.5557
$ws.Range("K113").Value = 2286.5557
$ws.Range("M113").Value = -116.5556999999999
$ws.Range("H122").Value = 2855
$ws.Range("I122").Value = 2433.875
$ws.Range("J122").Value = 3697.25
$ws.Range("K122").Value = 7301.625
$ws.Range("L122").Value = 11091.75
$ws.Range("M122").Value = -4851.625
$ws.Range("N122").Value = -15991.75
$ws.Range("H136").Value = 3283.1956
$ws.Range("I136").Value = 3188.361
$ws.Range("J136").Value = 3624.6
$ws.Range("K136").Value = 9565.082999999999
$ws.Range("L136").Value = 10873.8
$ws.Range("M136").Value = -7015.082999999999
$ws.Range("N136").Value = -15973.8
$ws.Range("H141").Value = 52500
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 60000
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -70360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 40008500
$ws.Range("J25").Value = 40008500
$ws.Range("L25").Value = 40008500
$ws.Range("N25").Value = -40009086
$ws.Range("H96").Value = 6688.091
$ws.Range("I96").Value = 8279.799999999999
$ws.Range("K96").Value = 8279.799999999999
$ws.Range("M96").Value = -6906.799999999999
